$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.23"
$ws.Range("E2").Value = "'-0.41%"
$ws.Range("G2").Value = "'18"
$ws.Range("D3").Value = "'39.04"
$ws.Range("E3").Value = "'7.32%"
$ws.Range("G3").Value = "'18"
$ws.Range("D4").Value = "'5.107"
$ws.Range("E4").Value = "'1.04%"
$ws.Range("G4").Value = "'18"
$ws.Range("D5").Value = "'0.08052"
$ws.Range("E5").Value = "'-0.42%"
$ws.Range("G5").Value = "'18"
$ws.Range("E6").Value = "'-7.69%"
$ws.Range("G6").Value = "'18"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.195"
$ws.Range("E7").Value = "'1.05%"
$ws.Range("G7").Value = "'18"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.002"
$ws.Range("E8").Value = "'2.03%"
$ws.Range("G8").Value = "'18"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9309"
$ws.Range("E9").Value = "'-0.06%"
$ws.Range("G9").Value = "'18"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1452"
$ws.Range("E10").Value = "'-1.58%"
$ws.Range("G10").Value = "'18"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1926"
$ws.Range("E11").Value = "'0.04%"
$ws.Range("G11").Value = "'18"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09101"
$ws.Range("E12").Value = "'-0.08%"
$ws.Range("G12").Value = "'18"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03501"
$ws.Range("E13").Value = "'1.38%"
$ws.Range("G13").Value = "'18"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09781"
$ws.Range("E14").Value = "'-1.32%"
$ws.Range("G14").Value = "'18"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001394"
$ws.Range("E15").Value = "'-1.11%"
$ws.Range("G15").Value = "'18"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005801"
$ws.Range("E16").Value = "'-8.12%"
$ws.Range("G16").Value = "'18"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.785"
$ws.Range("E17").Value = "'-1.64%"
$ws.Range("G17").Value = "'18"
$ws.Range("D18").Value = "'3.446"
$ws.Range("E18").Value = "'-1.15%"
$ws.Range("G18").Value = "'18"
$ws.Range("E19").Value = "'-1.00%"
$ws.Range("G19").Value = "'18"
$ws.Range("D20").Value = "'0.1302"
$ws.Range("E20").Value = "'-0.35%"
$ws.Range("G20").Value = "'18"
$ws.Range("D21").Value = "'4.787"
$ws.Range("E21").Value = "'-0.62%"
$ws.Range("G21").Value = "'18"
$ws.Range("D22").Value = "'0.2413"
$ws.Range("E22").Value = "'3.11%"
$ws.Range("G22").Value = "'18"
$ws.Range("D23").Value = "'0.04383"
$ws.Range("E23").Value = "'0.52%"
$ws.Range("G23").Value = "'18"
$ws.Range("E24").Value = "'0.43%"
$ws.Range("G24").Value = "'18"
$ws.Range("D25").Value = "'0.004280"
$ws.Range("E25").Value = "'-12.99%"
$ws.Range("G25").Value = "'18"
$ws.Range("E26").Value = "'0.14%"
$ws.Range("G26").Value = "'18"
$ws.Range("G27").Value = "'18"
$ws.Range("G28").Value = "'18"
$ws.Range("G29").Value = "'18"
$ws.Range("G30").Value = "'18"
$ws.Range("G31").Value = "'18"
$ws.Range("G32").Value = "'18"
$ws.Range("G33").Value = "'18"
$ws.Range("G34").Value = "'18"
$ws.Range("G35").Value = "'18"
$ws.Range("G36").Value = "'18"
$ws.Range("G37").Value = "'18"
$ws.Range("G38").Value = "'18"
$ws.Range("D39").Value = "'0.02039"
$ws.Range("E39").Value = "'1.25%"
$ws.Range("G39").Value = "'18"
$ws.Range("D40").Value = "'0.05033"
$ws.Range("E40").Value = "'-3.07%"
$ws.Range("G40").Value = "'18"
$ws.Range("D41").Value = "'0.007440"
$ws.Range("E41").Value = "'-0.76%"
$ws.Range("G41").Value = "'18"
$ws.Range("D42").Value = "'0.01021"
$ws.Range("E42").Value = "'0.78%"
$ws.Range("G42").Value = "'18"
$ws.Range("D43").Value = "'0.1348"
$ws.Range("E43").Value = "'-1.49%"
$ws.Range("G43").Value = "'18"
$ws.Range("E44").Value = "'-2.17%"
$ws.Range("G44").Value = "'18"
$ws.Range("D45").Value = "'0.009028"
$ws.Range("E45").Value = "'-9.36%"
$ws.Range("G45").Value = "'18"
$ws.Range("D46").Value = "'0.00006205"
$ws.Range("E46").Value = "'-1.23%"
$ws.Range("G46").Value = "'18"
$ws.Range("E47").Value = "'0.13%"
$ws.Range("G47").Value = "'18"
$ws.Range("D48").Value = "'0.002976"
$ws.Range("G48").Value = "'18"
$ws.Range("E49").Value = "'28.19%"
$ws.Range("G49").Value = "'18"
$ws.Range("E50").Value = "'0.13%"
$ws.Range("G50").Value = "'18"
$ws.Range("E51").Value = "'0.13%"
$ws.Range("G51").Value = "'18"
